$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column G, row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 03:01:03"

# zh-cn sheet: "Correspond Handoff Datetime" column H, row 2
#              "Correspond Handback DateTime" column K, row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-19 03:00:56"
$wsZhCn.Range("K2").Value = "2016-08-19 03:01:28"

# de-de sheet: "Correspond Handback DateTime" column K, row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-19 03:01:35"
